$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00002074986032285508
$ws.Range("C2").Value = 0.04240448674262143
$ws.Range("D2").Value = 3.900430680208489
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 649.2701327467715
